$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '55.557.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.951.69'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -5.26%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.60%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '489.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.945.11'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.418'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.16'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.100'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.345'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -8.88%  '
$ws.Range("E13").Value = '  +1.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.470.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '24.64'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '55.637.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.967.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000139'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.61'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '316.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -8.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.460'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '60.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -11.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.162'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("E28").Value = '  +0.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0835'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.45'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.50'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.41%  '
$ws.Range("E32").Value = '  -2.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.69'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.39'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -9.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '149.67'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.35'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -9.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.29'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.65'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0646'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.90'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -10.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.991.24'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '36.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.52%  '
$ws.Range("E44").Value = '  -3.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.629'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.121.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0233'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.06'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -9.62%  '
